$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")
$ws.Activate()
